# INTL 201 Exam 1, Part 2 essay - revise prompt paragraph formatting and
# replace the unfinished answer with the full essay text (adds a 3rd paragraph).
#
# Each paragraph's exact WordprocessingML (runs + paragraph/run formatting) is
# applied via Range.InsertXML so the run-by-run formatting from the source
# matches precisely; this is the standard Word COM way to push fully-formed
# OOXML into a Range.

$d = $word.ActiveDocument

# --- Paragraph 1 (the essay prompt/question) ---
# Drop the "ListParagraph" numbering + "Strong" character style/shading, and
# instead make the run bold + centered, sized 12pt (sz/szCs = 24 half-points).
$para1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>How did the Core and the Periphery shift over time? Consider the history of our connected world. Provide explicit details examples of the shift from your readings and lecture notes. In your answer make sure to explain the Core and Periphery relationship, as proposed by Wallerstein. Please also describe the factors that contributed to the shift, according to the class lectures and course readings.</w:t></w:r></w:p>'
$d.Paragraphs.Item(1).Range.InsertXML($para1Xml)

# --- Paragraph 2 (the answer) ---
# Replace the placeholder "The perpsad" text with the full first answer
# paragraph (firstLine indent instead of the old list indent).
$para2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Core and Periphery</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> are</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> not separate entities – they operate simultaneously and emulat</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>e</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> a </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>symbiotic relationship</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. Although this dynamic is mostly defined, </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">countries and powers will shift from being a Periphery to a Core and vice versa. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>As said by Dennis (2023), a</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ccording to the Core and Periphery relationship proposed by Wallerstein,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the Core is responsible for manufacturing high value goods, while the Periphery is responsible for producing raw materials, as well as exporting them. Thus, the Core exports their manufactured goods to the Periphery and the Periphery export their raw materials to the Core (Lecture). It can be seen that the Core is more powerful than the Periphery, despite their reliance on one another. In the present day it is obvious how America is a Core due to their high value manufactured exports, but when America was only the 13 colonies, they were a Periphery. How this shift happens is due to a variety of factors. Reasons, such as trade routes, new technologies, and abundance </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">or lack </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">of raw materials </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>drove the shifts between the Core and the Periphery.</w:t></w:r></w:p>'
$d.Paragraphs.Item(2).Range.InsertXML($para2Xml)

# --- Paragraph 3 (new) ---
# Add a new paragraph after paragraph 2 continuing the essay.
$d.Paragraphs.Item(2).Range.InsertParagraphAfter()
$para3Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>To analyze how the Core and the Periphery shifted over time, let us look back to the 1</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>00s</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (11</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Century). Dennis (2023) reiterated how Southern Spain shifted from a Periphery into a Core. The main reasons for this shift was due to</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$d.Paragraphs.Item(3).Range.InsertXML($para3Xml)
